$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.344.65"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.048.08"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.10"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.66"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.70"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").Value = "2.331.82"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.64"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "2.044.19"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "37.231.32"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.39"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  -4.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.73"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.21"
$ws.Range("E28").Value = "  -6.71%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0614"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  -4.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.25"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.98"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.89"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.475.89"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0942"
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.37"
$ws.Range("E45").Value = "  -5.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.93"
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.11"
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("D51").Value = "2.232.86"
$ws.Range("E51").Value = "  -1.62%  "
